$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.181.86"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.637.51"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'216.83"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.516"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'20.07"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "1.866.29"
$ws.Range("D13").Value = "1.647.08"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'0.541"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'66.58"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "27.166.00"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'216.88"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'6.83"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'9.12"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "'147.43"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "'3.03"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "1.301.75"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.550"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +5.72%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "1.776.55"
$ws.Range("D45").Value = "'62.16"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'91.24"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'7.64"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -0.80%  "
